$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the Year column (C2:C5) from 2017 to 2019
$ws.Range("C2").Value = 2019
$ws.Range("C3").Value = 2019
$ws.Range("C4").Value = 2019
$ws.Range("C5").Value = 2019

# Update the active cell / selection to C6
$ws.Range("C6").Select()
